# Rebuild the template paragraphs for "Haptics Showroom":
#  - two new empty bold paragraphs at the top
#  - the (formerly Heading1) title paragraph now carries body indentation and
#    hosts its drawing/rectangle whose text box contains the actual heading
#    text ("Haptics Showroom") instead of a loose run after the drawing
#  - a third blank paragraph and a second copy of the rectangle+heading
#    text box (carrying the _GoBack bookmark) replace the old bookmark-only
#    paragraph at the end of the document
$d = $word.ActiveDocument

$xml = @'
<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="de-DE"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251661312" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="58993D77" wp14:editId="09D7ED34"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>40640</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>116840</wp:posOffset></wp:positionV><wp:extent cx="5623560" cy="1151890"/><wp:effectExtent l="0" t="0" r="15240" b="10160"/><wp:wrapNone/><wp:docPr id="2" name="Rectangle 2"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5623560" cy="1151890"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="240"/><w:jc w:val="center"/><w:rPr><w:sz w:val="96"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="96"/></w:rPr><w:t>Haptics Showroom</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect id="Rectangle 2" o:spid="_x0000_s1026" style="position:absolute;left:0;text-align:left;margin-left:3.2pt;margin-top:9.2pt;width:442.8pt;height:90.7pt;z-index:251661312;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQDGVvA/egIAAEUFAAAOAAAAZHJzL2Uyb0RvYy54bWysVFFP2zAQfp+0/2D5faTJKGMVKapATJMQ&#xA;Q8DEs3HsJpLt885u0+7X7+ykAQHaw7Q+uD7f3Xe+L9/57HxnDdsqDB24mpdHM86Uk9B0bl3znw9X&#xA;n045C1G4RhhwquZ7Ffj58uOHs94vVAUtmEYhIxAXFr2veRujXxRFkK2yIhyBV46cGtCKSCauiwZF&#xA;T+jWFNVsdlL0gI1HkCoEOr0cnHyZ8bVWMv7QOqjITM3pbjGvmNentBbLM7FYo/BtJ8driH+4hRWd&#xA;o6IT1KWIgm2wewNlO4kQQMcjCbYArTupcg/UTTl71c19K7zKvRA5wU80hf8HK2+2t8i6hr4dZ05Y&#xA;+kR3RJpwa6NYmejpfVhQ1L2/xdEKtE297jTa9E9dsF2mdD9RqnaRSTqcn1SfT6s5Z5J8ZTmvqnkm&#xA;vXhO9xjiNwWWpU3NkcpnKsX2OkQqSaGHkFTNwVVnTDpPNxvukndxb1QKMO5OaWqJqlcZKItJXRhk&#xA;W0EyEFIqF8vB1YpGDcfzGf1Sw1RvyshWBkzImgpP2CNAEupb7AFmjE+pKmtxSp797WJD8pSRK4OL&#xA;U7LtHOB7AIa6GisP8QeSBmoSS0/Q7OmDIwyTELy86oj2axHirUCSPg0JjXP8QYs20Nccxh1nLeDv&#xA;985TPCmSvJz1NEo1D782AhVn5rsjrX4tj4/T7GXjeP6lIgNfep5eetzGXgB9JtIj3S5vU3w0h61G&#xA;sI809atUlVzCSapdcxnxYFzEYcTp3ZBqtcphNG9exGt372UCT6wmWT3sHgX6UXuRZHsDh7ETi1cS&#xA;HGJTpoPVJoLusj6feR35plnNwhnflfQYvLRz1PPrt/wDAAD//wMAUEsDBBQABgAIAAAAIQBkNPLv&#xA;3gAAAAkBAAAPAAAAZHJzL2Rvd25yZXYueG1sTE/LTsMwELwj8Q/WVuJGnRYUtSFOVSpxAiqlASRu&#xA;rr0kgXgdxW4b+vUsJzitZmc0j3w1uk4ccQitJwWzaQICyXjbUq3gpXq4XoAIUZPVnSdU8I0BVsXl&#xA;Ra4z609U4nEXa8EmFDKtoImxz6QMpkGnw9T3SMx9+MHpyHCopR30ic1dJ+dJkkqnW+KERve4adB8&#xA;7Q5OAb6+fZbn90ezfTJrX9ImVvfVs1JXk3F9ByLiGP/E8Fufq0PBnfb+QDaIjnGSspLvLS9gfrFM&#xA;5yD2/FimNyCLXP5fUPwAAAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAA&#xA;CwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAxlbwP3oCAABFBQAA&#xA;DgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYACAAAACEAZDTy794AAAAJ&#xA;AQAADwAAAAAAAAAAAAAAAADUBAAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAAN8FAAAA&#xA;AA==&#xA;" filled="f" strokecolor="#243f60 [1604]" strokeweight="2pt"><v:textbox><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="240"/><w:jc w:val="center"/><w:rPr><w:sz w:val="96"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="96"/></w:rPr><w:t>Haptics Showroom</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p/><w:p/><w:p/><w:p><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="de-DE"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251663360" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3669AAFE" wp14:editId="0ED19C3B"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>37465</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>494665</wp:posOffset></wp:positionV><wp:extent cx="5623560" cy="1151890"/><wp:effectExtent l="0" t="0" r="15240" b="10160"/><wp:wrapNone/><wp:docPr id="4" name="Rectangle 4"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5623560" cy="1151890"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="240"/><w:jc w:val="center"/><w:rPr><w:sz w:val="96"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:sz w:val="96"/></w:rPr><w:t>Haptics Showroom</w:t></w:r></w:p><w:bookmarkEnd w:id="0"/><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect id="Rectangle 4" o:spid="_x0000_s1027" style="position:absolute;margin-left:2.95pt;margin-top:38.95pt;width:442.8pt;height:90.7pt;z-index:251663360;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQDGVvA/egIAAEUFAAAOAAAAZHJzL2Uyb0RvYy54bWysVFFP2zAQfp+0/2D5faTJKGMVKapATJMQ&#xA;Q8DEs3HsJpLt885u0+7X7+ykAQHaw7Q+uD7f3Xe+L9/57HxnDdsqDB24mpdHM86Uk9B0bl3znw9X&#xA;n045C1G4RhhwquZ7Ffj58uOHs94vVAUtmEYhIxAXFr2veRujXxRFkK2yIhyBV46cGtCKSCauiwZF&#xA;T+jWFNVsdlL0gI1HkCoEOr0cnHyZ8bVWMv7QOqjITM3pbjGvmNentBbLM7FYo/BtJ8driH+4hRWd&#xA;o6IT1KWIgm2wewNlO4kQQMcjCbYArTupcg/UTTl71c19K7zKvRA5wU80hf8HK2+2t8i6hr4dZ05Y&#xA;+kR3RJpwa6NYmejpfVhQ1L2/xdEKtE297jTa9E9dsF2mdD9RqnaRSTqcn1SfT6s5Z5J8ZTmvqnkm&#xA;vXhO9xjiNwWWpU3NkcpnKsX2OkQqSaGHkFTNwVVnTDpPNxvukndxb1QKMO5OaWqJqlcZKItJXRhk&#xA;W0EyEFIqF8vB1YpGDcfzGf1Sw1RvyshWBkzImgpP2CNAEupb7AFmjE+pKmtxSp797WJD8pSRK4OL&#xA;U7LtHOB7AIa6GisP8QeSBmoSS0/Q7OmDIwyTELy86oj2axHirUCSPg0JjXP8QYs20Nccxh1nLeDv&#xA;985TPCmSvJz1NEo1D782AhVn5rsjrX4tj4/T7GXjeP6lIgNfep5eetzGXgB9JtIj3S5vU3w0h61G&#xA;sI809atUlVzCSapdcxnxYFzEYcTp3ZBqtcphNG9exGt372UCT6wmWT3sHgX6UXuRZHsDh7ETi1cS&#xA;HGJTpoPVJoLusj6feR35plnNwhnflfQYvLRz1PPrt/wDAAD//wMAUEsDBBQABgAIAAAAIQBkNPLv&#xA;3gAAAAkBAAAPAAAAZHJzL2Rvd25yZXYueG1sTE/LTsMwELwj8Q/WVuJGnRYUtSFOVSpxAiqlASRu&#xA;rr0kgXgdxW4b+vUsJzitZmc0j3w1uk4ccQitJwWzaQICyXjbUq3gpXq4XoAIUZPVnSdU8I0BVsXl&#xA;Ra4z609U4nEXa8EmFDKtoImxz6QMpkGnw9T3SMx9+MHpyHCopR30ic1dJ+dJkkqnW+KERve4adB8&#xA;7Q5OAb6+fZbn90ezfTJrX9ImVvfVs1JXk3F9ByLiGP/E8Fufq0PBnfb+QDaIjnGSspLvLS9gfrFM&#xA;5yD2/FimNyCLXP5fUPwAAAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAA&#xA;AAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAA&#xA;CwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAxlbwP3oCAABFBQAA&#xA;DgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYACAAAACEAZDTy794AAAAJ&#xA;AQAADwAAAAAAAAAAAAAAAADUBAAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAAN8FAAAA&#xA;AA==&#xA;" filled="f" strokecolor="#243f60 [1604]" strokeweight="2pt"><v:textbox><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:spacing w:before="240"/><w:jc w:val="center"/><w:rPr><w:sz w:val="96"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:r><w:rPr><w:sz w:val="96"/></w:rPr><w:t>Haptics Showroom</w:t></w:r></w:p><w:bookmarkEnd w:id="1"/><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p>
'@

$d.Content.InsertXML($xml)
